# Insert a new row at position 189 (shifts existing rows 189-259 down to 190-260)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(189).Insert()

# Populate the newly inserted row 189 with the new record's data
$ws.Cells.Item(189, 1).Value = 10
$ws.Cells.Item(189, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(189, 3).Value = "La Araucanía"
$ws.Cells.Item(189, 4).Value = 44784
$ws.Cells.Item(189, 5).Value = 9
$ws.Cells.Item(189, 6).Value = 100112052
$ws.Cells.Item(189, 7).Value = "Albahaca"
$ws.Cells.Item(189, 8).Value = "Sin especificar"
$ws.Cells.Item(189, 9).Value = "Primera"
$ws.Cells.Item(189, 10).Value = 100
$ws.Cells.Item(189, 11).Value = 6000
$ws.Cells.Item(189, 12).Value = 6000
$ws.Cells.Item(189, 13).Value = 6000
$ws.Cells.Item(189, 14).Value = '$/paquete'
$ws.Cells.Item(189, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(189, 16).Value = 6000
$ws.Cells.Item(189, 17).Value = 1
$ws.Cells.Item(189, 18).Value = "Hortaliza"
